# Insert a new data row at row 375 (pushing the existing rows 375-407 down
# to 376-408) and populate it with the new weekly price observation.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows.Item(375).Insert()

$ws.Cells.Item(375, 1).Value = 4
$ws.Cells.Item(375, 2).Value = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(375, 3).Value = "Los Lagos"
$ws.Cells.Item(375, 4).Value = 44769
$ws.Cells.Item(375, 5).Value = 10
$ws.Cells.Item(375, 6).Value = 100114013
$ws.Cells.Item(375, 7).Value = "Zanahoria"
$ws.Cells.Item(375, 8).Value = "Sin especificar"
$ws.Cells.Item(375, 9).Value = "Primera"
$ws.Cells.Item(375, 10).Value = 80
$ws.Cells.Item(375, 11).Value = 10000
$ws.Cells.Item(375, 12).Value = 10000
$ws.Cells.Item(375, 13).Value = 10000
$ws.Cells.Item(375, 14).Value = "$/saco 20 kilos"
$ws.Cells.Item(375, 15).Value = "Provincia de Llanquihue"
$ws.Cells.Item(375, 16).Value = 500
$ws.Cells.Item(375, 17).Value = 20
$ws.Cells.Item(375, 18).Value = "Hortaliza"
